# Add a new "rest" row (C2010, 2025-05-15 ~ 2025-05-30) to the
# machine_rest sheet (3rd sheet in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# New data row under the existing header row.
$ws.Range("A2").Value = "C2010"
$ws.Range("B2").Value = 45792
$ws.Range("C2").Value = 45807

# Date-format the two new date cells (maps to built-in numFmtId 14).
# Format B2 first, then copy/paste its format onto C2 so both cells
# share a single cell-style record (rather than each getting its own).
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Copy()
[void]$ws.Range("C2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen columns B:C to fit the new date values/headers.
$ws.Columns("B:C").ColumnWidth = 10.428571428571429

# Leave the selection where the author left it when saving.
[void]$ws.Range("B3").Select()
